$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Every Status cell (col F) that still reads "In Document" becomes "Tentative",
#        picking up the same text/formatting already used for "Tentative" elsewhere
#        in the column (e.g. F2) so it matches the other green "Tentative" cells. ---
$tentativeRef = $ws.Range("F2")
$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "In Document") {
        $cell.Value2 = $tentativeRef.Value2
        $cell.Font.Color = $tentativeRef.Font.Color
    }
}

# --- 2. Filter the table (A1:F37) down to the "Team Member" (col E) = "Jack" rows. ---
$ws.Range("A1:F37").AutoFilter(5, @("Jack"), 7)

# --- 3. Leave the selection on F28, a visible "Jack" row once the filter is applied. ---
$ws.Range("F28").Select()
